$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# row 58
$ws.Range("H58").Value = 5104270
$ws.Range("I58").Value = 71428570
$ws.Range("J58").Value = 2400.6924
$ws.Range("K58").Value = 214285710
$ws.Range("L58").Value = 7202.0772
$ws.Range("M58").Value = -214285560
$ws.Range("N58").Value = -7502.0772
# row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("L87").ClearContents()
# row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("N90").Value = 0
$ws.Range("L90").ClearContents()
# row 93
$ws.Range("H93").Value = 49800
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 49800
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 49800
$ws.Range("N93").Value = -54792
# row 98
$ws.Range("H98").Value = 45089.777
$ws.Range("I98").Value = 1444.9166
$ws.Range("J98").Value = 132379.5
$ws.Range("K98").Value = 1444.9166
$ws.Range("L98").Value = 132379.5
$ws.Range("M98").Value = 53.08339999999998
$ws.Range("N98").Value = -135375.5
# row 106
$ws.Range("H106").Value = 100169.625
$ws.Range("I106").Value = 2995
$ws.Range("J106").Value = 197344.25
$ws.Range("K106").Value = 2995
$ws.Range("L106").Value = 197344.25
$ws.Range("M106").Value = -2364
$ws.Range("N106").Value = -198606.25
# row 122
$ws.Range("H122").Value = 45089.777
$ws.Range("I122").Value = 1444.9166
$ws.Range("J122").Value = 132379.5
$ws.Range("K122").Value = 4334.7498
$ws.Range("L122").Value = 397138.5
$ws.Range("M122").Value = -1884.7498
$ws.Range("N122").Value = -402038.5
# row 132
$ws.Range("H132").Value = 29772.156
$ws.Range("I132").Value = 4278.8335
$ws.Range("J132").Value = 106252.125
$ws.Range("K132").Value = 12836.5005
$ws.Range("L132").Value = 318756.375
$ws.Range("M132").Value = -10306.5005
$ws.Range("N132").Value = -323816.375
# row 141
$ws.Range("H141").Value = 5349.231
$ws.Range("I141").Value = 3334.2856
$ws.Range("J141").Value = 7700
$ws.Range("K141").Value = 10002.8568
$ws.Range("L141").Value = 23100
$ws.Range("M141").Value = -4822.856800000001
$ws.Range("N141").Value = -33460

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# row 4
$ws.Range("H4").Value = 800
$ws.Range("I4").Value = 800
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 800
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -684
# row 23
$ws.Range("H23").Value = 15000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 15000
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = 15000
$ws.Range("N23").Value = -15518
$ws.Range("L23").ClearContents()
# row 32
$ws.Range("H32").Value = 11225.246
$ws.Range("I32").Value = 10061.393
$ws.Range("J32").Value = 21118
$ws.Range("K32").Value = 10061.393
$ws.Range("L32").Value = 21118
$ws.Range("M32").Value = -9774.393
$ws.Range("N32").Value = -21692
# row 37
$ws.Range("H37").Value = 33578.332
$ws.Range("I37").Value = 15000
$ws.Range("J37").Value = 42867.5
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 42867.5
$ws.Range("M37").Value = -14727
$ws.Range("N37").Value = -43413.5
# row 44
$ws.Range("H44").Value = 29505
$ws.Range("I44").Value = 5000
$ws.Range("J44").Value = 37673.332
$ws.Range("K44").Value = 5000
$ws.Range("L44").Value = 37673.332
$ws.Range("M44").Value = -4512
$ws.Range("N44").Value = -38649.332
# row 55
$ws.Range("H55").Value = 37000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 37000
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = 37000
$ws.Range("N55").Value = -37630
$ws.Range("L55").ClearContents()
# row 61
$ws.Range("H61").Value = 2822.6128
$ws.Range("I61").Value = 2123
$ws.Range("J61").Value = 3398.7646
$ws.Range("K61").Value = 2123
$ws.Range("L61").Value = 3398.7646
$ws.Range("M61").Value = -1911
$ws.Range("N61").Value = -3822.7646
# row 63
$ws.Range("H63").Value = 2506.84
$ws.Range("I63").Value = 2061.476
$ws.Range("J63").Value = 4845
$ws.Range("K63").Value = 2061.476
$ws.Range("L63").Value = 4845
$ws.Range("M63").Value = -1375.476
$ws.Range("N63").Value = -6217
# row 66
$ws.Range("H66").Value = 2506.84
$ws.Range("I66").Value = 2061.476
$ws.Range("J66").Value = 4845
$ws.Range("K66").Value = 10307.38
$ws.Range("L66").Value = 24225
$ws.Range("M66").Value = -6875.380000000001
$ws.Range("N66").Value = -31089
# row 80
$ws.Range("H80").Value = 47863
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 47863
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 47863
$ws.Range("N80").Value = -49859
# row 83
$ws.Range("H83").Value = 47863
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 47863
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 143589
$ws.Range("N83").Value = -153573
# row 102
$ws.Range("H102").Value = 34254.21
$ws.Range("I102").Value = 5603.3335
$ws.Range("J102").Value = 39626.25
$ws.Range("K102").Value = 5603.3335
$ws.Range("L102").Value = 39626.25
$ws.Range("M102").Value = -3981.3335
$ws.Range("N102").Value = -42870.25
# row 103
$ws.Range("H103").Value = 30000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 30000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
# row 122
$ws.Range("H122").Value = 1658.2
$ws.Range("I122").Value = 1543.409
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4630.227000000001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2180.227000000001
$ws.Range("N122").Value = -12400
# row 132
$ws.Range("H132").Value = 11113422
$ws.Range("I132").Value = 16668108
$ws.Range("J132").Value = 4050.6667
$ws.Range("K132").Value = 50004324
$ws.Range("L132").Value = 12152.0001
$ws.Range("M132").Value = -50001794
$ws.Range("N132").Value = -17212.0001
# row 136
$ws.Range("H136").Value = 2822.6128
$ws.Range("I136").Value = 2123
$ws.Range("J136").Value = 3398.7646
$ws.Range("K136").Value = 6369
$ws.Range("L136").Value = 10196.2938
$ws.Range("M136").Value = -3819
$ws.Range("N136").Value = -15296.2938

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
# row 15
$ws.Range("H15").Value = 39125
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 39125
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 39125
$ws.Range("N15").Value = -39579
# row 19
$ws.Range("H19").Value = 22785.715
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 22785.715
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 22785.715
$ws.Range("N19").Value = -23131.715
# row 34
$ws.Range("H34").Value = 2485
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2485
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 2485
$ws.Range("N34").Value = -2713
# row 35
$ws.Range("H35").Value = 30764.445
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 30764.445
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 30764.445
$ws.Range("N35").Value = -31384.445
# row 40
$ws.Range("H40").Value = 18655.555
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 18655.555
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 18655.555
$ws.Range("N40").Value = -19185.555

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# row 6
$ws.Range("H6").Value = 4000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 4000
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = 4000
$ws.Range("N6").Value = -4226
$ws.Range("L6").ClearContents()
# row 7
$ws.Range("H7").Value = 501.57144
$ws.Range("I7").Value = 583.5
$ws.Range("J7").Value = 10
$ws.Range("K7").Value = 583.5
$ws.Range("L7").Value = 10
$ws.Range("M7").Value = -470.5
$ws.Range("N7").Value = -236
# row 17
$ws.Range("H17").Value = 34995
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 34995
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = 34995
$ws.Range("N17").Value = -35343
$ws.Range("L17").ClearContents()
# row 25
$ws.Range("H25").Value = 35217.555
$ws.Range("I25").Value = 30000
$ws.Range("J25").Value = 39391.6
$ws.Range("K25").Value = 30000
$ws.Range("L25").Value = 39391.6
$ws.Range("M25").Value = -29826
$ws.Range("N25").Value = -39739.6
# row 41
$ws.Range("H41").Value = 27171.334
$ws.Range("I41").Value = 18059
$ws.Range("J41").Value = 28993.8
$ws.Range("K41").Value = 18059
$ws.Range("L41").Value = 28993.8
$ws.Range("M41").Value = -17631
$ws.Range("N41").Value = -29849.8
# row 50
$ws.Range("H50").Value = 39187.25
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 39187.25
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 39187.25
$ws.Range("N50").Value = -40437.25
# row 51
$ws.Range("H51").Value = 83364200
$ws.Range("I51").Value = 250015040
$ws.Range("J51").Value = 38779.75
$ws.Range("K51").Value = 250015040
$ws.Range("L51").Value = 38779.75
$ws.Range("M51").Value = -250014304
$ws.Range("N51").Value = -40251.75
# row 58
$ws.Range("H58").Value = 1854.921
$ws.Range("I58").Value = 1148.6154
$ws.Range("J58").Value = 3385.25
$ws.Range("K58").Value = 1148.6154
$ws.Range("L58").Value = 3385.25
$ws.Range("M58").Value = -945.6153999999999
$ws.Range("N58").Value = -3791.25
# row 59
$ws.Range("H59").Value = 42639.668
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 42639.668
$ws.Range("K59").Value = 0
$ws.Range("M59").Value = 42639.668
$ws.Range("N59").Value = -44929.668
$ws.Range("L59").ClearContents()
# row 60
$ws.Range("H60").Value = 11103
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 11103
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = 11103
$ws.Range("N60").Value = -12125
$ws.Range("L60").ClearContents()
# row 61
$ws.Range("H61").Value = 83364200
$ws.Range("I61").Value = 250015040
$ws.Range("J61").Value = 38779.75
$ws.Range("K61").Value = 250015040
$ws.Range("L61").Value = 38779.75
$ws.Range("M61").Value = -250014692
$ws.Range("N61").Value = -39475.75
# row 68
$ws.Range("H68").Value = 30000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 30000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
# row 71
$ws.Range("H71").Value = 30000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 30000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
# row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("L74").ClearContents()
# row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("L77").ClearContents()
# row 136
$ws.Range("H136").Value = 1854.921
$ws.Range("I136").Value = 1148.6154
$ws.Range("J136").Value = 3385.25
$ws.Range("K136").Value = 3445.8462
$ws.Range("L136").Value = 10155.75
$ws.Range("M136").Value = -895.8462
$ws.Range("N136").Value = -15255.75

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# row 17
$ws.Range("H17").Value = 280
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 333.33334
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 1000.00002
$ws.Range("M17").Value = -431
$ws.Range("N17").Value = -1338.00002
# row 34
$ws.Range("H34").Value = 1833.2727
$ws.Range("I34").Value = 700
$ws.Range("J34").Value = 1946.6
$ws.Range("K34").Value = 2100
$ws.Range("L34").Value = 5839.799999999999
$ws.Range("M34").Value = -2016
$ws.Range("N34").Value = -6007.799999999999
# row 39
$ws.Range("H39").Value = 571.4286
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 571.4286
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 1714.2858
$ws.Range("N39").Value = -2302.2858
# row 55
$ws.Range("H55").Value = 2214.8572
$ws.Range("I55").Value = 704
$ws.Range("J55").Value = 2466.6667
$ws.Range("K55").Value = 2112
$ws.Range("L55").Value = 7400.000100000001
$ws.Range("M55").Value = -1935
$ws.Range("N55").Value = -7754.000100000001
# row 101
$ws.Range("H101").Value = 7857.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 7857.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 23572.5
$ws.Range("N101").Value = -28440.5

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# row 18
$ws.Range("H18").Value = 11666.667
$ws.Range("I18").Value = 10000
# row 43
$ws.Range("H43").Value = 13017
$ws.Range("I43").Value = 13017
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 13017
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = -12866
$ws.Range("M43").ClearContents()
# row 46
$ws.Range("H46").Value = 23403.947
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 24426.389
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 24426.389
$ws.Range("M46").Value = -4844
$ws.Range("N46").Value = -24738.389
# row 57
$ws.Range("H57").Value = 13400.952
$ws.Range("I57").Value = 6318.1816
$ws.Range("J57").Value = 21192
$ws.Range("K57").Value = 6318.1816
$ws.Range("L57").Value = 21192
$ws.Range("M57").Value = -5498.1816
$ws.Range("N57").Value = -22832
# row 80
$ws.Range("H80").Value = 6421.3687
$ws.Range("I80").Value = 4692.3076
$ws.Range("J80").Value = 10167.667
$ws.Range("K80").Value = 4692.3076
$ws.Range("L80").Value = 10167.667
$ws.Range("M80").Value = -3694.3076
$ws.Range("N80").Value = -12163.667
# row 83
$ws.Range("H83").Value = 6421.3687
$ws.Range("I83").Value = 4692.3076
$ws.Range("J83").Value = 10167.667
$ws.Range("K83").Value = 23461.538
$ws.Range("L83").Value = 50838.335
$ws.Range("M83").Value = -18469.538
$ws.Range("N83").Value = -60822.335
# row 97
$ws.Range("H97").Value = 2967.276
$ws.Range("I97").Value = 2075.5557
$ws.Range("J97").Value = 15005.5
$ws.Range("K97").Value = 2075.5557
$ws.Range("L97").Value = 15005.5
$ws.Range("M97").Value = -1579.5557
$ws.Range("N97").Value = -15997.5
# row 102
$ws.Range("H102").Value = 1325
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 1150
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 1150
$ws.Range("M102").Value = 122
$ws.Range("N102").Value = -4394
# row 122
$ws.Range("H122").Value = 1418
$ws.Range("I122").Value = 1428.2858
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 4284.857400000001
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -1834.857400000001
$ws.Range("N122").Value = -9100
# row 126
$ws.Range("H126").Value = 5700.4443
$ws.Range("I126").Value = 10237.333
$ws.Range("J126").Value = 2070.9333
$ws.Range("K126").Value = 30711.999
$ws.Range("L126").Value = 6212.7999
$ws.Range("M126").Value = -28241.999
$ws.Range("N126").Value = -11152.7999
# row 132
$ws.Range("H132").Value = 29414638
$ws.Range("I132").Value = 41668856
$ws.Range("J132").Value = 4519.3
$ws.Range("K132").Value = 125006568
$ws.Range("L132").Value = 13557.9
$ws.Range("M132").Value = -125004038
$ws.Range("N132").Value = -18617.9

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 2633.1177
$ws.Range("I7").Value = 1786.8
$ws.Range("J7").Value = 3842.1428
$ws.Range("K7").Value = 1786.8
$ws.Range("L7").Value = 3842.1428
$ws.Range("M7").Value = -1674.8
$ws.Range("N7").Value = -4066.1428
# row 18
$ws.Range("H18").Value = 1500
$ws.Range("I18").Value = 1500
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1500
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1328
# row 36
$ws.Range("H36").Value = 48715
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 48715
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 48715
$ws.Range("N36").Value = -49839
# row 40
$ws.Range("H40").Value = 3762.75
$ws.Range("I40").Value = 3762.75
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3762.75
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3626.75
# row 122
$ws.Range("H122").Value = 113410.664
$ws.Range("I122").Value = 145027.86
$ws.Range("J122").Value = 2750.5
$ws.Range("K122").Value = 435083.58
$ws.Range("L122").Value = 8251.5
$ws.Range("M122").Value = -432633.58
$ws.Range("N122").Value = -13151.5
# row 126
$ws.Range("H126").Value = 2633.1177
$ws.Range("I126").Value = 1786.8
$ws.Range("J126").Value = 3842.1428
$ws.Range("K126").Value = 5360.4
$ws.Range("L126").Value = 11526.4284
$ws.Range("M126").Value = -2890.4
$ws.Range("N126").Value = -16466.4284
# row 132
$ws.Range("H132").Value = 4230.5386
$ws.Range("I132").Value = 3300.1177
$ws.Range("J132").Value = 5988
$ws.Range("K132").Value = 9900.3531
$ws.Range("L132").Value = 17964
$ws.Range("M132").Value = -7370.3531
$ws.Range("N132").Value = -23024

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
# row 54
$ws.Range("H54").Value = 50000784
$ws.Range("I54").Value = 50000784
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 50000784
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = -50000264
$ws.Range("M54").ClearContents()
# row 81
$ws.Range("H81").Value = 1333.3334
$ws.Range("I81").Value = 1000
# row 84
$ws.Range("H84").Value = 1333.3334
$ws.Range("I84").Value = 1000
# row 122
$ws.Range("H122").Value = 3177558
$ws.Range("I122").Value = 3574552.8
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 10723658.4
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -10721208.4
$ws.Range("N122").Value = -9700
